# Applies the cryptos price/volume/coin-name updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($range, $value) {
    # Force the cell to stay text (e.g. "316.95", "1.76%") instead of Excel
    # auto-converting numeric/percent-looking strings to numbers, then restore
    # the default "Normal" style so no stray number-format/style id is left behind.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell 'D2' '316.95'
Set-TextCell 'E2' '1.76%'
Set-TextCell 'D3' '37.71'
Set-TextCell 'E3' '0.50%'
Set-TextCell 'D4' '5.164'
Set-TextCell 'E4' '0.60%'
Set-TextCell 'D5' '0.07966'
Set-TextCell 'E5' '1.44%'
Set-TextCell 'B6' 'GateToken'
Set-TextCell 'C6' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell 'D6' '4.457'
Set-TextCell 'E6' '0.68%'
Set-TextCell 'B7' 'KuCoinToken'
Set-TextCell 'C7' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextCell 'D7' '8.502'
Set-TextCell 'E7' '2.92%'
Set-TextCell 'B8' 'FTXToken'
Set-TextCell 'C8' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell 'D8' '1.939'
Set-TextCell 'E8' '1.96%'
Set-TextCell 'B9' 'BTSEToken'
Set-TextCell 'C9' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell 'D9' '2.977'
Set-TextCell 'E9' '2.37%'
Set-TextCell 'B10' 'MXToken'
Set-TextCell 'C10' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D10' '0.9401'
Set-TextCell 'E10' '2.32%'
Set-TextCell 'B11' 'LiechtensteinCryptoassetsExchange'
Set-TextCell 'C11' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 'D11' '0.1263'
Set-TextCell 'E11' '6.72%'
Set-TextCell 'B12' 'WazirX'
Set-TextCell 'C12' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 'D12' '0.1930'
Set-TextCell 'E12' '0.85%'
Set-TextCell 'B13' 'MandalaExchangeToken'
Set-TextCell 'C13' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 'D13' '0.09007'
Set-TextCell 'E13' '-0.97%'
Set-TextCell 'B14' 'BitrueCoin'
Set-TextCell 'C14' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 'D14' '0.03391'
Set-TextCell 'E14' '1.51%'
Set-TextCell 'B15' 'BitMartToken'
Set-TextCell 'C15' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 'D15' '0.09535'
Set-TextCell 'E15' '-0.52%'
Set-TextCell 'B16' 'BitForexToken'
Set-TextCell 'C16' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 'D16' '0.001369'
Set-TextCell 'E16' '-1.28%'
Set-TextCell 'B17' 'TigerCash'
Set-TextCell 'C17' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell 'D17' '0.006096'
Set-TextCell 'E17' '6.70%'
Set-TextCell 'B18' 'LEO'
Set-TextCell 'C18' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 'D18' '3.410'
Set-TextCell 'E18' '-2.90%'
Set-TextCell 'D19' '0.3513'
Set-TextCell 'E19' '2.07%'
Set-TextCell 'D20' '6.517'
Set-TextCell 'E20' '23.93%'
Set-TextCell 'E21' '2.64%'
Set-TextCell 'D22' '0.2300'
Set-TextCell 'E22' '-11.21%'
Set-TextCell 'D23' '0.04346'
Set-TextCell 'E23' '-0.52%'
Set-TextCell 'D24' '0.001198'
Set-TextCell 'E24' '-4.23%'
Set-TextCell 'D25' '0.004410'
Set-TextCell 'E25' '-5.92%'
Set-TextCell 'D26' '0.0001322'
Set-TextCell 'E26' '-3.15%'
Set-TextCell 'D27' '0.0003972'
Set-TextCell 'E27' '-0.43%'
Set-TextCell 'D39' '0.02354'
Set-TextCell 'E39' '1.99%'
Set-TextCell 'D40' '0.05172'
Set-TextCell 'E40' '2.02%'
Set-TextCell 'D41' '0.007407'
Set-TextCell 'E41' '-1.07%'
Set-TextCell 'D42' '0.1396'
Set-TextCell 'E42' '3.12%'
Set-TextCell 'D43' '0.008562'
Set-TextCell 'E43' '-5.35%'
Set-TextCell 'D44' '0.001991'
Set-TextCell 'E44' '1.68%'
Set-TextCell 'D45' '0.008750'
Set-TextCell 'E45' '-7.23%'
Set-TextCell 'D46' '0.00006407'
Set-TextCell 'E46' '-3.32%'
Set-TextCell 'D47' '0.00000000747'
Set-TextCell 'E47' '-0.82%'
Set-TextCell 'D48' '0.002849'
Set-TextCell 'E48' '-13.31%'
Set-TextCell 'D49' '0.001682'
Set-TextCell 'E49' '68.22%'
Set-TextCell 'D50' '0.00002090'
Set-TextCell 'E50' '-0.82%'
Set-TextCell 'D51' '0.0001991'
Set-TextCell 'E51' '-0.82%'
